# Auto-generated Excel COM-interop script
# Applies scheduled market-price/profit refresh values to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 391.90475
$ws.Range("I9").Value = 302.66666
$ws.Range("J9").Value = 615
$ws.Range("K9").Value = 302.66666
$ws.Range("L9").Value = 615
$ws.Range("M9").Value = -133.66666
$ws.Range("N9").Value = -953
$ws.Range("H17").Value = 35391.516
$ws.Range("J17").Value = 35391.516
$ws.Range("L17").Value = 106174.548
$ws.Range("N17").Value = -106510.548
$ws.Range("H34").Value = 1492.1428
$ws.Range("I34").Value = 1492.1428
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1492.1428
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -1289.1428
$ws.Range("H36").Value = 1492.1428
$ws.Range("I36").Value = 1492.1428
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1492.1428
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -777.1428000000001
$ws.Range("H100").Value = 2919.2104
$ws.Range("I100").Value = 2300.3333
$ws.Range("J100").Value = 5240
$ws.Range("K100").Value = 2300.3333
$ws.Range("L100").Value = 5240
$ws.Range("M100").Value = -1759.3333
$ws.Range("N100").Value = -6322
$ws.Range("H127").Value = 930
$ws.Range("I127").Value = 506.66666
$ws.Range("J127").Value = 2200
$ws.Range("K127").Value = 1519.99998
$ws.Range("L127").Value = 6600
$ws.Range("M127").Value = 3440.00002
$ws.Range("N127").Value = -16520
$ws.Range("H132").Value = 8003975
$ws.Range("I132").Value = 8337057.5
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 25011172.5
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -25008642.5
$ws.Range("N132").Value = -35060
$ws.Range("H138").Value = 4337.1665
$ws.Range("I138").Value = 1408.7391
$ws.Range("J138").Value = 5903.5347
$ws.Range("K138").Value = 4226.2173
$ws.Range("L138").Value = 17710.6041
$ws.Range("M138").Value = 913.7826999999997
$ws.Range("N138").Value = -27990.6041

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1441.97
$ws.Range("I32").Value = 1268.9333
$ws.Range("J32").Value = 2999.3
$ws.Range("K32").Value = 1268.9333
$ws.Range("L32").Value = 2999.3
$ws.Range("M32").Value = -981.9332999999999
$ws.Range("N32").Value = -3573.3
$ws.Range("H74").Value = 915.3570999999999
$ws.Range("I74").Value = 821.5
$ws.Range("K74").Value = 821.5
$ws.Range("M74").Value = 52.5
$ws.Range("H77").Value = 915.3570999999999
$ws.Range("I77").Value = 821.5
$ws.Range("K77").Value = 4107.5
$ws.Range("M77").Value = 260.5
$ws.Range("H101").Value = 39970.6
$ws.Range("J101").Value = 39970.6
$ws.Range("L101").Value = 39970.6
$ws.Range("N101").Value = -46460.6
$ws.Range("H122").Value = 2112.5625
$ws.Range("I122").Value = 1516.75
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 4550.25
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -2100.25
$ws.Range("N122").Value = -16600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 26169
$ws.Range("I15").Value = 3000
$ws.Range("J15").Value = 49338
$ws.Range("K15").Value = 3000
$ws.Range("L15").Value = 49338
$ws.Range("M15").Value = -2773
$ws.Range("N15").Value = -49792
$ws.Range("H134").Value = 4272.923
$ws.Range("I134").Value = 3304.8572
$ws.Range("J134").Value = 5402.3335
$ws.Range("K134").Value = 9914.571599999999
$ws.Range("L134").Value = 16207.0005
$ws.Range("M134").Value = -7379.571599999999
$ws.Range("N134").Value = -21277.0005
$ws.Range("H141").Value = 35030.668
$ws.Range("J141").Value = 29587.273
$ws.Range("L141").Value = 29587.273
$ws.Range("N141").Value = -39947.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 6000
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H22").Value = 1661.5834
$ws.Range("I22").Value = 490
$ws.Range("J22").Value = 1768.091
$ws.Range("K22").Value = 490
$ws.Range("L22").Value = 1768.091
$ws.Range("M22").Value = -140
$ws.Range("N22").Value = -2468.091
$ws.Range("H107").Value = 1732.4166
$ws.Range("I107").Value = 949
$ws.Range("J107").Value = 2292
$ws.Range("K107").Value = 949
$ws.Range("L107").Value = 2292
$ws.Range("M107").Value = 971
$ws.Range("N107").Value = -6132
$ws.Range("H132").Value = 2679.7568
$ws.Range("I132").Value = 2111.5356
$ws.Range("J132").Value = 4447.5557
$ws.Range("K132").Value = 6334.6068
$ws.Range("L132").Value = 13342.6671
$ws.Range("M132").Value = -3804.6068
$ws.Range("N132").Value = -18402.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2378.25
$ws.Range("I75").Value = 1013
$ws.Range("K75").Value = 3039
$ws.Range("M75").Value = -2041
$ws.Range("H78").Value = 2378.25
$ws.Range("I78").Value = 1013
$ws.Range("K78").Value = 9117
$ws.Range("M78").Value = -4125
$ws.Range("H87").Value = 14214.286
$ws.Range("H90").Value = 14214.286
$ws.Range("H120").Value = 15765
$ws.Range("H122").Value = 1433.762
$ws.Range("I122").Value = 400.2
$ws.Range("J122").Value = 2373.3635
$ws.Range("K122").Value = 3601.8
$ws.Range("L122").Value = 21360.2715
$ws.Range("M122").Value = -1151.8
$ws.Range("N122").Value = -26260.2715
$ws.Range("H131").Value = 1319.3939
$ws.Range("I131").Value = 1896.25
$ws.Range("J131").Value = 1134.8
$ws.Range("K131").Value = 5688.75
$ws.Range("L131").Value = 3404.4
$ws.Range("M131").Value = -648.75
$ws.Range("N131").Value = -13484.4
$ws.Range("H138").Value = 2196.7778
$ws.Range("I138").Value = 1311.8334
$ws.Range("J138").Value = 3966.6667
$ws.Range("K138").Value = 3935.5002
$ws.Range("L138").Value = 11900.0001
$ws.Range("M138").Value = 1204.4998
$ws.Range("N138").Value = -22180.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2956.7693
$ws.Range("I126").Value = 1979.5385
$ws.Range("J126").Value = 3934
$ws.Range("K126").Value = 5938.6155
$ws.Range("L126").Value = 11802
$ws.Range("M126").Value = -3468.6155
$ws.Range("N126").Value = -16742

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 50836.668
$ws.Range("J3").Value = 50836.668
$ws.Range("L3").Value = 50836.668
$ws.Range("N3").Value = -51060.668
$ws.Range("H15").Value = 50836.668
$ws.Range("J15").Value = 50836.668
$ws.Range("L15").Value = 50836.668
$ws.Range("N15").Value = -51176.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10310.951
$ws.Range("I132").Value = 2038.9811
$ws.Range("J132").Value = 59023.668
$ws.Range("K132").Value = 6116.9433
$ws.Range("L132").Value = 177071.004
$ws.Range("M132").Value = -3586.9433
$ws.Range("N132").Value = -182131.004
